$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C4/C5 previously held raw numeric values (400 / 500), with C4 carrying a
# custom numeric format. Replace them with the concatenated text values,
# matching the pattern already used by C1:C3 ("<A>.0<B>.0"), and drop the
# leftover custom formatting so the cells fall back to the default style.
$ws.Range("C4").Value = "400.040.0"
$ws.Range("C5").Value = "500.050.0"
$ws.Range("C4").Style = "Normal"
$ws.Range("C5").Style = "Normal"

# Select the full first row (A1:XFD1), matching the saved selection state.
$ws.Rows("1:1").Select()
